$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finish out row 7: add X7 (numeric) and Y7 (string "Down") ---
$ws.Range("X7").Value = -3.5699769999999944
$ws.Range("Y7").Value = "Down"

# --- Append new row 8 with a full record ---
$ws.Range("A8").Value = 42649.890902777777
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"

$ws.Range("B8").Value = -3
$ws.Range("C8").Value = "Neutral"
$ws.Range("D8").Value = 14
$ws.Range("E8").Value = 20255
$ws.Range("F8").Value = 3786
$ws.Range("G8").Value = 59
$ws.Range("H8").Value = 36
$ws.Range("I8").Value = 84
$ws.Range("J8").Value = 13
$ws.Range("K8").Value = 32916
$ws.Range("L8").Value = 456
$ws.Range("M8").Value = 278
$ws.Range("N8").Value = 184
$ws.Range("O8").Value = 30
$ws.Range("P8").Value = "Noun"

$ws.Range("Q8").Value = 49.72799223503381
$ws.Range("R8").Value = 0

$ws.Range("S8").Value = 0.1095
$ws.Range("S8").NumberFormat = "0.00%"

$ws.Range("T8").Value = 0.0025
$ws.Range("T8").NumberFormat = "0.00%"

$ws.Range("U8").Value = 5.95
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0
